# Refresh the scraped crypto price/volume figures (Tue Feb 20 04:19:09 UTC 2024 run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.717.90'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '2.918.64'
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''353.87'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").Value = '''110.32'
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("D7").Value = '''0.571'
$ws.Range("E7").Value = '  +1.92%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("D10").Value = '''39.14'
$ws.Range("E10").Value = '  -1.98%  '
$ws.Range("E11").Value = '  +2.96%  '
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '''19.63'
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("D14").Value = '''7.87'
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").Value = '3.376.78'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").Value = '2.918.06'
$ws.Range("E16").Value = '  +1.58%  '
$ws.Range("D17").Value = '''0.975'
$ws.Range("E17").Value = '  -1.75%  '
$ws.Range("D18").Value = '51.691.38'
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").Value = '''7.51'
$ws.Range("E19").Value = '  -2.65%  '
$ws.Range("D20").Value = '''3.25'
$ws.Range("E20").Value = '  -2.67%  '
$ws.Range("D21").Value = '''13.91'
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '''70.67'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("D24").Value = '''269.45'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '''2.80'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").Value = '''0.182'
$ws.Range("E26").Value = '  +11.96%  '
$ws.Range("D27").Value = '''26.98'
$ws.Range("E27").Value = '  +2.40%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("D29").Value = '''7.28'
$ws.Range("E29").Value = '  +12.47%  '
$ws.Range("D30").Value = '''0.105'
$ws.Range("E30").Value = '  +11.19%  '
$ws.Range("D31").Value = '''10.51'
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").Value = '''38.66'
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("E33").Value = '  -2.01%  '
$ws.Range("D34").Value = '''52.10'
$ws.Range("E34").Value = '  -2.14%  '
$ws.Range("E35").Value = '  -4.45%  '
$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '''1.89'
$ws.Range("E37").Value = '  -16.10%  '
$ws.Range("E38").Value = '  -2.92%  '
$ws.Range("D39").Value = '''18.29'
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("D40").Value = '''2.00'
$ws.Range("E40").Value = '  -1.79%  '
$ws.Range("D41").Value = '''2.74'
$ws.Range("E41").Value = '  +4.28%  '
$ws.Range("D42").Value = '''0.120'
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("D43").Value = '''22.69'
$ws.Range("E43").Value = '  +0.92%  '
$ws.Range("D44").Value = '''120.12'
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("D47").Value = '''3.43'
$ws.Range("E47").Value = '  -3.70%  '
$ws.Range("D48").Value = '2.130.35'
$ws.Range("E48").Value = '  -3.04%  '
$ws.Range("D49").Value = '''0.247'
$ws.Range("E49").Value = '  -7.84%  '
$ws.Range("D50").Value = '''0.0329'
$ws.Range("E50").Value = '  +3.20%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '''60.43'
$ws.Range("E51").Value = '  +2.86%  '
